$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new timesheet entry on row 11
$ws.Range("A11").Value = "Wk[12] Saturday 2.6.18"
$ws.Range("B11").Value = "1300  - 1600"
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "Reparing logic of updating the facotry"

# Update the active selection to D6
$ws.Range("D6").Select() | Out-Null
